$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.71895933333333
$ws.Range("H2").Value = 50.156878
$ws.Range("I2").Value = 0.02912144738161902
$ws.Range("J2").Value = 0.03059269312988411
$ws.Range("M2").Value = 279.6078796666666
$ws.Range("N2").Value = 838.823639
$ws.Range("O2").Value = 0.5726675140320879
$ws.Range("P2").Value = 0.5775260287976519
$ws.Range("Q2").Value = 4674.75276942656
$ws.Range("R2").Value = 42072.77492483904
$ws.Range("S2").Value = 0.01667690687704802
$ws.Range("T2").Value = 0.01766807657352718
$ws.Range("G3").Value = 16.71895933333333
$ws.Range("H3").Value = 50.156878
$ws.Range("I3").Value = 0.02912144738161902
$ws.Range("J3").Value = 0.03059269312988411
$ws.Range("O3").Value = 0.001411539217074187
$ws.Range("P3").Value = 0.001423514724607417
$ws.Range("Q3").Value = 11.52256187488578
$ws.Range("R3").Value = 103.703056873972
$ws.Range("S3").Value = 0.00004110606503711763
$ws.Range("T3").Value = 0.0000435491491357862
$ws.Range("G4").Value = 16.71895933333333
$ws.Range("H4").Value = 50.156878
$ws.Range("I4").Value = 0.02912144738161902
$ws.Range("J4").Value = 0.03059269312988411
$ws.Range("M4").Value = 55.12872433333333
$ws.Range("N4").Value = 165.386173
$ws.Range("O4").Value = 0.1129096560274583
$ws.Range("P4").Value = 0.1138675822543568
$ws.Range("Q4").Value = 921.6949002275438
$ws.Range("R4").Value = 8295.254102047893
$ws.Range("S4").Value = 0.003288092606880328
$ws.Range("T4").Value = 0.003483516001349377
$ws.Range("G5").Value = 16.71895933333333
$ws.Range("H5").Value = 50.156878
$ws.Range("I5").Value = 0.02912144738161902
$ws.Range("J5").Value = 0.03059269312988411
$ws.Range("M5").Value = 12.3225355
$ws.Range("N5").Value = 24.645071
$ws.Range("O5").Value = 0.02523790023288966
$ws.Range("P5").Value = 0.01696801249072354
$ws.Range("Q5").Value = 206.0199699080563
$ws.Range("R5").Value = 1236.119819448338
$ws.Range("S5").Value = 0.0007349641836546467
$ws.Range("T5").Value = 0.0005190971991527459
$ws.Range("G6").Value = 16.71895933333333
$ws.Range("H6").Value = 50.156878
$ws.Range("I6").Value = 0.02912144738161902
$ws.Range("J6").Value = 0.03059269312988411
$ws.Range("M6").Value = 140.5068483333333
$ws.Range("N6").Value = 421.520545
$ws.Range("O6").Value = 0.2877733904904901
$ws.Range("P6").Value = 0.2902148617326603
$ws.Range("Q6").Value = 2349.128283339835
$ws.Range("R6").Value = 21142.15455005851
$ws.Range("S6").Value = 0.008380377648998908
$ws.Range("T6").Value = 0.008878454206719023
$ws.Range("I7").Value = 0.2708539632042961
$ws.Range("J7").Value = 0.2845377865576845
$ws.Range("M7").Value = 279.6078796666666
$ws.Range("N7").Value = 838.823639
$ws.Range("O7").Value = 0.5726675140320879
$ws.Range("P7").Value = 0.5775260287976519
$ws.Range("Q7").Value = 43479.13405563185
$ws.Range("R7").Value = 391312.2065006867
$ws.Range("S7").Value = 0.1551092657739429
$ws.Range("T7").Value = 0.1643279779135335
$ws.Range("I8").Value = 0.2708539632042961
$ws.Range("J8").Value = 0.2845377865576845
$ws.Range("O8").Value = 0.001411539217074187
$ws.Range("P8").Value = 0.001423514724607417
$ws.Range("S8").Value = 0.0003823209911628326
$ws.Range("T8").Value = 0.0004050437288720663
$ws.Range("I9").Value = 0.2708539632042961
$ws.Range("J9").Value = 0.2845377865576845
$ws.Range("M9").Value = 55.12872433333333
$ws.Range("N9").Value = 165.386173
$ws.Range("O9").Value = 0.1129096560274583
$ws.Range("P9").Value = 0.1138675822543568
$ws.Range("Q9").Value = 8572.538078907099
$ws.Range("R9").Value = 77152.84271016388
$ws.Range("S9").Value = 0.0305820278190709
$ws.Range("T9").Value = 0.03239962981532978
$ws.Range("I10").Value = 0.2708539632042961
$ws.Range("J10").Value = 0.2845377865576845
$ws.Range("M10").Value = 12.3225355
$ws.Range("N10").Value = 24.645071
$ws.Range("O10").Value = 0.02523790023288966
$ws.Range("P10").Value = 0.01696801249072354
$ws.Range("Q10").Value = 1916.159063716309
$ws.Range("R10").Value = 11496.95438229786
$ws.Range("S10").Value = 0.006835785301032792
$ws.Range("T10").Value = 0.00482804071639362
$ws.Range("I11").Value = 0.2708539632042961
$ws.Range("J11").Value = 0.2845377865576845
$ws.Range("M11").Value = 140.5068483333333
$ws.Range("N11").Value = 421.520545
$ws.Range("O11").Value = 0.2877733904904901
$ws.Range("P11").Value = 0.2902148617326603
$ws.Range("Q11").Value = 21848.86957299734
$ws.Range("R11").Value = 196639.826156976
$ws.Range("S11").Value = 0.07794456331908671
$ws.Range("T11").Value = 0.08257709438355561
$ws.Range("G12").Value = 194.8548433333333
$ws.Range("H12").Value = 584.56453
$ws.Range("I12").Value = 0.3394024086099587
$ws.Range("J12").Value = 0.3565493705749576
$ws.Range("M12").Value = 279.6078796666666
$ws.Range("N12").Value = 838.823639
$ws.Range("O12").Value = 0.5726675140320879
$ws.Range("P12").Value = 0.5775260287976519
$ws.Range("Q12").Value = 54482.94958721384
$ws.Range("R12").Value = 490346.5462849246
$ws.Range("S12").Value = 0.194364733595168
$ws.Range("T12").Value = 0.2059165420584576
$ws.Range("G13").Value = 194.8548433333333
$ws.Range("H13").Value = 584.56453
$ws.Range("I13").Value = 0.3394024086099587
$ws.Range("J13").Value = 0.3565493705749576
$ws.Range("O13").Value = 0.001411539217074187
$ws.Range("P13").Value = 0.001423514724607417
$ws.Range("Q13").Value = 134.2922692833578
$ws.Range("R13").Value = 1208.63042355022
$ws.Range("S13").Value = 0.0004790798101223944
$ws.Range("T13").Value = 0.0005075532790629586
$ws.Range("G14").Value = 194.8548433333333
$ws.Range("H14").Value = 584.56453
$ws.Range("I14").Value = 0.3394024086099587
$ws.Range("J14").Value = 0.3565493705749576
$ws.Range("M14").Value = 55.12872433333333
$ws.Range("N14").Value = 165.386173
$ws.Range("O14").Value = 0.1129096560274583
$ws.Range("P14").Value = 0.1138675822543568
$ws.Range("Q14").Value = 10742.09894313819
$ws.Range("R14").Value = 96678.89048824368
$ws.Range("S14").Value = 0.03832180921104128
$ws.Range("T14").Value = 0.04059941478168314
$ws.Range("G15").Value = 194.8548433333333
$ws.Range("H15").Value = 584.56453
$ws.Range("I15").Value = 0.3394024086099587
$ws.Range("J15").Value = 0.3565493705749576
$ws.Range("M15").Value = 12.3225355
$ws.Range("N15").Value = 24.645071
$ws.Range("O15").Value = 0.02523790023288966
$ws.Range("P15").Value = 0.01696801249072354
$ws.Range("Q15").Value = 2401.105724321938
$ws.Range("R15").Value = 14406.63434593163
$ws.Range("S15").Value = 0.00856580412730059
$ws.Range("T15").Value = 0.006049934173475497
$ws.Range("G16").Value = 194.8548433333333
$ws.Range("H16").Value = 584.56453
$ws.Range("I16").Value = 0.3394024086099587
$ws.Range("J16").Value = 0.3565493705749576
$ws.Range("M16").Value = 140.5068483333333
$ws.Range("N16").Value = 421.520545
$ws.Range("O16").Value = 0.2877733904904901
$ws.Range("P16").Value = 0.2902148617326603
$ws.Range("Q16").Value = 27378.43991925209
$ws.Range("R16").Value = 246405.9592732688
$ws.Range("S16").Value = 0.09767098186632651
$ws.Range("T16").Value = 0.1034759262822783
$ws.Range("G17").Value = 82.82950199999999
$ws.Range("H17").Value = 165.659004
$ws.Range("I17").Value = 0.1442742299952585
$ws.Range("J17").Value = 0.1010420758958371
$ws.Range("M17").Value = 279.6078796666666
$ws.Range("N17").Value = 838.823639
$ws.Range("O17").Value = 0.5726675140320879
$ws.Range("P17").Value = 0.5775260287976519
$ws.Range("Q17").Value = 23159.78142806592
$ws.Range("R17").Value = 138958.6885683955
$ws.Range("S17").Value = 0.08262116463027837
$ws.Range("T17").Value = 0.05835442883359377
$ws.Range("G18").Value = 82.82950199999999
$ws.Range("H18").Value = 165.659004
$ws.Range("I18").Value = 0.1442742299952585
$ws.Range("J18").Value = 0.1010420758958371
$ws.Range("O18").Value = 0.001411539217074187
$ws.Range("P18").Value = 0.001423514724607417
$ws.Range("Q18").Value = 57.085374922716
$ws.Range("R18").Value = 342.5122495362959
$ws.Range("S18").Value = 0.0002036487336514883
$ws.Range("T18").Value = 0.0001438348828426243
$ws.Range("G19").Value = 82.82950199999999
$ws.Range("H19").Value = 165.659004
$ws.Range("I19").Value = 0.1442742299952585
$ws.Range("J19").Value = 0.1010420758958371
$ws.Range("M19").Value = 55.12872433333333
$ws.Range("N19").Value = 165.386173
$ws.Range("O19").Value = 0.1129096560274583
$ws.Range("P19").Value = 0.1138675822543568
$ws.Range("Q19").Value = 4566.284782425281
$ws.Range("R19").Value = 27397.70869455169
$ws.Range("S19").Value = 0.01628995368239104
$ws.Range("T19").Value = 0.0115054168882202
$ws.Range("G20").Value = 82.82950199999999
$ws.Range("H20").Value = 165.659004
$ws.Range("I20").Value = 0.1442742299952585
$ws.Range("J20").Value = 0.1010420758958371
$ws.Range("M20").Value = 12.3225355
$ws.Range("N20").Value = 24.645071
$ws.Range("O20").Value = 0.02523790023288966
$ws.Range("P20").Value = 0.01696801249072354
$ws.Range("Q20").Value = 1020.669478842321
$ws.Range("R20").Value = 4082.677915369283
$ws.Range("S20").Value = 0.003641178622797311
$ws.Range("T20").Value = 0.001714483205889201
$ws.Range("G21").Value = 82.82950199999999
$ws.Range("H21").Value = 165.659004
$ws.Range("I21").Value = 0.1442742299952585
$ws.Range("J21").Value = 0.1010420758958371
$ws.Range("M21").Value = 140.5068483333333
$ws.Range("N21").Value = 421.520545
$ws.Range("O21").Value = 0.2877733904904901
$ws.Range("P21").Value = 0.2902148617326603
$ws.Range("Q21").Value = 11638.11227503953
$ws.Range("R21").Value = 69828.67365023716
$ws.Range("S21").Value = 0.0415182843261403
$ws.Range("T21").Value = 0.02932391208529134
$ws.Range("G22").Value = 124.2078576666667
$ws.Range("H22").Value = 372.623573
$ws.Range("I22").Value = 0.2163479508088675
$ws.Range("J22").Value = 0.2272780738416368
$ws.Range("M22").Value = 279.6078796666666
$ws.Range("N22").Value = 838.823639
$ws.Range("O22").Value = 0.5726675140320879
$ws.Range("P22").Value = 0.5775260287976519
$ws.Range("Q22").Value = 34729.49572011579
$ws.Range("R22").Value = 312565.4614810421
$ws.Range("S22").Value = 0.1238954431556506
$ws.Range("T22").Value = 0.13125900341854
$ws.Range("G23").Value = 124.2078576666667
$ws.Range("H23").Value = 372.623573
$ws.Range("I23").Value = 0.2163479508088675
$ws.Range("J23").Value = 0.2272780738416368
$ws.Range("O23").Value = 0.001411539217074187
$ws.Range("P23").Value = 0.001423514724607417
$ws.Range("Q23").Value = 85.60297903576689
$ws.Range("R23").Value = 770.4268113219019
$ws.Range("S23").Value = 0.0003053836171003536
$ws.Range("T23").Value = 0.0003235336846939819
$ws.Range("G24").Value = 124.2078576666667
$ws.Range("H24").Value = 372.623573
$ws.Range("I24").Value = 0.2163479508088675
$ws.Range("J24").Value = 0.2272780738416368
$ws.Range("M24").Value = 55.12872433333333
$ws.Range("N24").Value = 165.386173
$ws.Range("O24").Value = 0.1129096560274583
$ws.Range("P24").Value = 0.1138675822543568
$ws.Range("Q24").Value = 6847.420745339569
$ws.Range("R24").Value = 61626.78670805612
$ws.Range("S24").Value = 0.02442777270807469
$ws.Range("T24").Value = 0.02587960476777437
$ws.Range("G25").Value = 124.2078576666667
$ws.Range("H25").Value = 372.623573
$ws.Range("I25").Value = 0.2163479508088675
$ws.Range("J25").Value = 0.2272780738416368
$ws.Range("M25").Value = 12.3225355
$ws.Range("N25").Value = 24.645071
$ws.Range("O25").Value = 0.02523790023288966
$ws.Range("P25").Value = 0.01696801249072354
$ws.Range("Q25").Value = 1530.555735476447
$ws.Range("R25").Value = 9183.334412858681
$ws.Range("S25").Value = 0.00546016799810432
$ws.Range("T25").Value = 0.003856457195812481
$ws.Range("G26").Value = 124.2078576666667
$ws.Range("H26").Value = 372.623573
$ws.Range("I26").Value = 0.2163479508088675
$ws.Range("J26").Value = 0.2272780738416368
$ws.Range("M26").Value = 140.5068483333333
$ws.Range("N26").Value = 421.520545
$ws.Range("O26").Value = 0.2877733904904901
$ws.Range("P26").Value = 0.2902148617326603
$ws.Range("Q26").Value = 17452.05461897859
$ws.Range("R26").Value = 157068.4915708073
$ws.Range("S26").Value = 0.06225918332993757
$ws.Range("T26").Value = 0.06595947477481598
